$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mojibake for "Turkiye": the UTF-8 bytes of "Türkiye" were
# re-interpreted/re-encoded as Latin-1/CP1252, producing "TÃ¼rkiye".
# The country name lives in column B, on the row whose code (column A) is "TUR".
$ws.Range("B117").Value = "TÃ¼rkiye"

# The dataset had three accidental duplicate rows (same Codigo/Pais pair
# repeated back to back). Remove the duplicated copies, deleting from the
# bottom up so the row numbers of the earlier duplicates are unaffected.
$ws.Rows.Item(143).Delete()
$ws.Rows.Item(57).Delete()
$ws.Rows.Item(40).Delete()
